# Weekly data refresh: a new week's price observation is added at the top
# of the data block (row 219), pushing all the existing observations for
# rows 219-349 down by one row (they become rows 220-350).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 219..349 down by one row to make room for the new weekly entry.
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A219").Value = 9
$ws.Range("B219").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C219").Value = "Metropolitana"
$ws.Range("D219").Value = 44606
$ws.Range("E219").Value = 13
$ws.Range("F219").Value = 100112039
$ws.Range("G219").Value = "Ciboulette"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 106
$ws.Range("K219").Value = 1000
$ws.Range("L219").Value = 1200
$ws.Range("M219").Value = 1100
$ws.Range("N219").Value = "`$/docena de atados"
$ws.Range("O219").Value = "Región Metropolitana"
$ws.Range("P219").Value = 367
$ws.Range("Q219").Value = 3
$ws.Range("R219").Value = "Hortaliza"
